# Updates the "cryptos" price/volume snapshot sheet to refreshed values
# (Price column D, Volume(1h) column E), and swaps the FTXToken /
# TrustWalletToken rows (45 <-> 46) to match the new ranking order.
#
# Price/Volume cells store their text verbatim (e.g. "234.40", "37.862.36")
# rather than as numbers, so a leading apostrophe is used to force Excel to
# keep them as literal text instead of auto-coercing to a number (which
# would silently drop meaningful trailing zeros / reparse dotted values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.850.51"
$ws.Range("D3").Value = "'2.088.23"
$ws.Range("E3").Value = "'  +1.01%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'234.40"
$ws.Range("E5").Value = "'  -0.08%  "
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("D7").Value = "'59.18"
$ws.Range("E7").Value = "'  +3.13%  "
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("D9").Value = "'0.391"
$ws.Range("E9").Value = "'  -1.24%  "
$ws.Range("D10").Value = "'0.0794"
$ws.Range("E10").Value = "'  +2.61%  "
$ws.Range("E11").Value = "'  +3.09%  "
$ws.Range("D12").Value = "'2.397.93"
$ws.Range("E12").Value = "'  +1.04%  "
$ws.Range("D13").Value = "'14.72"
$ws.Range("E13").Value = "'  +2.11%  "
$ws.Range("D14").Value = "'21.46"
$ws.Range("E14").Value = "'  +3.66%  "
$ws.Range("D15").Value = "'0.771"
$ws.Range("E15").Value = "'  -0.69%  "
$ws.Range("E16").Value = "'  +2.53%  "
$ws.Range("D17").Value = "'2.088.45"
$ws.Range("E17").Value = "'  +0.95%  "
$ws.Range("D18").Value = "'37.773.67"
$ws.Range("E18").Value = "'  +1.18%  "
$ws.Range("D19").Value = "'6.27"
$ws.Range("E19").Value = "'  +0.07%  "
$ws.Range("D20").Value = "'71.74"
$ws.Range("E20").Value = "'  +3.04%  "
$ws.Range("E21").Value = "'  +1.61%  "
$ws.Range("D22").Value = "'228.92"
$ws.Range("E22").Value = "'  +0.96%  "
$ws.Range("E23").Value = "'  +0.04%  "
$ws.Range("E24").Value = "'  -0.68%  "
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "'  +0.01%  "
$ws.Range("D26").Value = "'170.81"
$ws.Range("E26").Value = "'  +2.28%  "
$ws.Range("E27").Value = "'  +9.76%  "
$ws.Range("D28").Value = "'9.07"
$ws.Range("E28").Value = "'  +2.50%  "
$ws.Range("E29").Value = "'  +0.19%  "
$ws.Range("D30").Value = "'19.59"
$ws.Range("E30").Value = "'  +2.49%  "
$ws.Range("E31").Value = "'  +2.13%  "
$ws.Range("E32").Value = "'  +3.92%  "
$ws.Range("D33").Value = "'0.0632"
$ws.Range("E33").Value = "'  +2.50%  "
$ws.Range("D34").Value = "'4.70"
$ws.Range("E34").Value = "'  +3.50%  "
$ws.Range("D35").Value = "'2.51"
$ws.Range("E35").Value = "'  +0.52%  "
$ws.Range("E36").Value = "'  +5.90%  "
$ws.Range("E37").Value = "'  +2.72%  "
$ws.Range("E38").Value = "'  -0.12%  "
$ws.Range("E39").Value = "'  -3.17%  "
$ws.Range("D40").Value = "'0.0990"
$ws.Range("E40").Value = "'  +2.68%  "
$ws.Range("D41").Value = "'99.53"
$ws.Range("E41").Value = "'  +1.62%  "
$ws.Range("E43").Value = "'  +1.62%  "
$ws.Range("D44").Value = "'1.463.72"
$ws.Range("E44").Value = "'  -1.21%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.18"
$ws.Range("E45").Value = "'  +1.08%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.30"
$ws.Range("E46").Value = "'  +6.63%  "
$ws.Range("D47").Value = "'16.26"
$ws.Range("E47").Value = "'  +6.28%  "
$ws.Range("E48").Value = "'  +5.41%  "
$ws.Range("E49").Value = "'  +3.51%  "
$ws.Range("E50").Value = "'  +2.94%  "
$ws.Range("D51").Value = "'47.75"
$ws.Range("E51").Value = "'  +6.46%  "
